$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.30860000000002
$ws.Range("D4").Value = -7.466100000000004
$ws.Range("E4").Value = 12.47450000000001

$ws.Range("D5").Value = -8.220399999999994

$ws.Range("A7").Value = -21.45640000000001

$ws.Range("D8").Value = -8.485499999999996

$ws.Range("E9").Value = 13.42070000000001

$ws.Range("A16").Value = -20.14799999999998
$ws.Range("D16").Value = -8.211900000000005

$ws.Range("E18").Value = 12.7636
